$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 5.5.2 indicator label lost its trailing period after "2"
$ws.Range("B4").Value = "5.5.2 Доля женщин на руководящих должностях"

# Organization website corrected to the new domain
$ws.Range("B10").Value = "www.stat.gov.kg"

# Restore the selection to B12 as left by the author before saving
$ws.Range("B12").Select()
